$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Udi")

$xlPasteFormats = -4122

for ($r = 2; $r -le 131; $r++) {
    $mCell = $ws.Cells.Item($r, 13)   # column M
    $sCell = $ws.Cells.Item($r, 19)   # column S

    # Copy the formatting of the M cell onto the S cell (covers the
    # style-index change for S, which always ends up matching M's style).
    $mCell.Copy()
    $sCell.PasteSpecial($xlPasteFormats)

    # Move the comment text itself from M to S.
    $sCell.Value2 = $mCell.Value2

    # Clear the source cell's content (keeping its formatting/style).
    $mCell.ClearContents()
}

$excel.CutCopyMode = $false

# Update the view state recorded in the sheet (scroll position + selection).
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("O3").Select()
